# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighboring "sum" header cell (G1) so the
# new header cell reuses the existing bold/border/center style instead of
# creating a brand-new style entry.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
